# Auto-generated edit script: Add data for 2023-11-27
# Updates 2023 (column J) YTD violent-crime figures across the Citywide Totals,
# By Neighborhood summary, and per-neighborhood sheets, plus a couple of minor
# 2018/2019 (E/F) reclassification corrections on the Citywide/By-Neighborhood
# and Andersonville sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("J2").Value = 6975
$ws.Range("J3").Value = 7376
$ws.Range("J4").Value = 1603
$ws.Range("E5").Value = 538
$ws.Range("F5").Value = 463
$ws.Range("J5").Value = 577
$ws.Range("J6").Value = 9955
$ws.Range("E7").Value = 23718
$ws.Range("F7").Value = 21908
$ws.Range("J7").Value = 26486

$ws = $wb.Worksheets.Item(11)
$ws.Range("J2").Value = 68
$ws.Range("J6").Value = 261
$ws.Range("J7").Value = 391

$ws = $wb.Worksheets.Item(12)
$ws.Range("J6").Value = 593
$ws.Range("J7").Value = 1658

$ws = $wb.Worksheets.Item(13)
$ws.Range("J3").Value = 199
$ws.Range("J7").Value = 529

$ws = $wb.Worksheets.Item(14)
$ws.Range("J3").Value = 396
$ws.Range("J6").Value = 423
$ws.Range("J7").Value = 1197

$ws = $wb.Worksheets.Item(15)
$ws.Range("J2").Value = 136
$ws.Range("J7").Value = 385

$ws = $wb.Worksheets.Item(16)
$ws.Range("J6").Value = 239
$ws.Range("J7").Value = 820

$ws = $wb.Worksheets.Item(17)
$ws.Range("J3").Value = 185
$ws.Range("J6").Value = 248
$ws.Range("J7").Value = 668

$ws = $wb.Worksheets.Item(18)
$ws.Range("J3").Value = 163
$ws.Range("J7").Value = 403

$ws = $wb.Worksheets.Item(2)
$ws.Range("J4").Value = 121
$ws.Range("J7").Value = 761
$ws.Range("J8").Value = 1658
$ws.Range("J11").Value = 467
$ws.Range("J13").Value = 32
$ws.Range("J15").Value = 323
$ws.Range("J17").Value = 35
$ws.Range("J18").Value = 219
$ws.Range("J19").Value = 768
$ws.Range("J21").Value = 70
$ws.Range("J29").Value = 1423
$ws.Range("J31").Value = 272
$ws.Range("J33").Value = 1197
$ws.Range("J37").Value = 820
$ws.Range("J40").Value = 57
$ws.Range("J41").Value = 194
$ws.Range("J42").Value = 1144
$ws.Range("J43").Value = 226
$ws.Range("E48").Value = 271
$ws.Range("F48").Value = 276
$ws.Range("J48").Value = 300
$ws.Range("J51").Value = 323
$ws.Range("J52").Value = 680
$ws.Range("J53").Value = 391
$ws.Range("J54").Value = 511
$ws.Range("J55").Value = 418
$ws.Range("J57").Value = 124
$ws.Range("J60").Value = 155
$ws.Range("J65").Value = 668
$ws.Range("J67").Value = 984
$ws.Range("J68").Value = 59
$ws.Range("J76").Value = 380
$ws.Range("J79").Value = 738
$ws.Range("J83").Value = 529
$ws.Range("J85").Value = 1089
$ws.Range("J88").Value = 284
$ws.Range("J89").Value = 332
$ws.Range("J91").Value = 304
$ws.Range("J94").Value = 289
$ws.Range("J95").Value = 385
$ws.Range("J96").Value = 286
$ws.Range("J97").Value = 242
$ws.Range("J98").Value = 198
$ws.Range("J99").Value = 403
$ws.Range("E101").Value = 23718
$ws.Range("F101").Value = 21908
$ws.Range("J101").Value = 26486

$ws = $wb.Worksheets.Item(20)
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item(21)
$ws.Range("J2").Value = 252
$ws.Range("J3").Value = 366
$ws.Range("J7").Value = 984

$ws = $wb.Worksheets.Item(24)
$ws.Range("J2").Value = 127
$ws.Range("J4").Value = 37
$ws.Range("J6").Value = 239
$ws.Range("J7").Value = 511

$ws = $wb.Worksheets.Item(25)
$ws.Range("J2").Value = 427
$ws.Range("J3").Value = 505
$ws.Range("J6").Value = 362
$ws.Range("J7").Value = 1423

$ws = $wb.Worksheets.Item(26)
$ws.Range("J4").Value = 48
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 2
$ws.Range("E7").Value = 271
$ws.Range("F7").Value = 276
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item(27)
$ws.Range("J2").Value = 186
$ws.Range("J3").Value = 221
$ws.Range("J7").Value = 768

$ws = $wb.Worksheets.Item(29)
$ws.Range("J2").Value = 66
$ws.Range("J3").Value = 80
$ws.Range("J4").Value = 29
$ws.Range("J7").Value = 380

$ws = $wb.Worksheets.Item(31)
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item(32)
$ws.Range("J4").Value = 48
$ws.Range("J6").Value = 608
$ws.Range("J7").Value = 1144

$ws = $wb.Worksheets.Item(33)
$ws.Range("J5").Value = 16
$ws.Range("J6").Value = 32

$ws = $wb.Worksheets.Item(36)
$ws.Range("J3").Value = 80
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 236
$ws.Range("J7").Value = 418

$ws = $wb.Worksheets.Item(4)
$ws.Range("J2").Value = 86
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 286

$ws = $wb.Worksheets.Item(40)
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item(41)
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item(42)
$ws.Range("J2").Value = 206
$ws.Range("J7").Value = 738

$ws = $wb.Worksheets.Item(45)
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item(46)
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item(5)
$ws.Range("J2").Value = 238
$ws.Range("J3").Value = 230
$ws.Range("J4").Value = 31
$ws.Range("J6").Value = 243
$ws.Range("J7").Value = 761

$ws = $wb.Worksheets.Item(51)
$ws.Range("J2").Value = 55
$ws.Range("J6").Value = 153
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item(54)
$ws.Range("J2").Value = 92
$ws.Range("J6").Value = 143
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item(55)
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item(6)
$ws.Range("J6").Value = 221
$ws.Range("J7").Value = 467

$ws = $wb.Worksheets.Item(65)
$ws.Range("J6").Value = 167
$ws.Range("J7").Value = 242

$ws = $wb.Worksheets.Item(68)
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 284

$ws = $wb.Worksheets.Item(7)
$ws.Range("J4").Value = 33
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 332

$ws = $wb.Worksheets.Item(75)
$ws.Range("J6").Value = 130
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item(76)
$ws.Range("J2").Value = 25
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item(77)
$ws.Range("J2").Value = 30
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item(78)
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item(79)
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item(8)
$ws.Range("J2").Value = 292
$ws.Range("J3").Value = 390
$ws.Range("J4").Value = 70
$ws.Range("J7").Value = 1089

$ws = $wb.Worksheets.Item(89)
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item(9)
$ws.Range("J2").Value = 159
$ws.Range("J6").Value = 292
$ws.Range("J7").Value = 680

$ws = $wb.Worksheets.Item(90)
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 121
